# Set origin to sprite center for flying objects
# Rewrites the two-column (offset, size) table on the sheet with new
# sprite-center-relative values, drops the now-unused last row, renames
# the sheet, and removes the "Zarez" (comma) number-format style that
# used to be applied to column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet1 -> List1
$ws.Name = "List1"

# New data (column A: offset, column B: size) - replaces the old
# formula-driven values with plain literals, and the table now has 9
# rows instead of 10.
$values = @(
    @(1, 3),
    @(3, 4),
    @(7, 2),
    @(13, 2),
    @(17, 2),
    @(22, 2),
    @(28, 2),
    @(32, 2),
    @(43, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# Drop the old 10th row entirely (new table is only 9 rows).
$ws.Rows.Item(10).Clear()

# Column B no longer uses the "Zarez" comma style - restore plain/default
# formatting on the cells that still carry it.
$ws.Range("B1:B9").ClearFormats()

# Remove the now-unused "Zarez" cell style definition from the workbook.
foreach ($s in $wb.Styles) {
    if ($s.Name -eq "Zarez") {
        $s.Delete()
    }
}

# Match the new selection left in the file (A2).
$ws.Range("A2").Select()
